$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to list many kanji/word rows (rows 1-16) each paired with a
# "meaning" comment in column B. Trim it down to just 4 rows and rewrite the
# text - a single kanji/word per row plus a "no comment" placeholder - and
# blank out what used to be rows 5-16 entirely.
$ws.Range("A1:B16").ClearContents()

$ws.Range("A1").Value = "郷"
$ws.Range("B1").Value = "Không có comment"
$ws.Range("A2").Value = "愁"
$ws.Range("B2").Value = "Không có comment"
$ws.Range("A3").Value = "幾"
$ws.Range("B3").Value = "Không có comment"
$ws.Range("A4").Value = "要"
$ws.Range("B4").Value = "Không có comment"

# Leave the selection on B4, matching the saved view's active cell.
$ws.Range("B4").Select()
